# Update the "startup" sheet query text per the commit:
# "updated keyword for case files tab data writing issue + more curated scripts"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- CasesTab row (row 2) ---
# B2: cases query gains a trailing `Cohort` column
$ws.Range("B2").Value = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE diag.primary_disease_site IN [''Mouth'']
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`,
coalesce(co.cohort_description, '''') AS `Cohort`'
# C2: StatQuery replaced with the new program/study/case/sample/file counter
$ws.Range("C2").Value = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.primary_disease_site IN [''Mouth'']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'

# --- SamplesTab row (row 3) ---
# B3 (samples query) text is unchanged.
# C3: StatQuery replaced with the new program/study/case/sample/file counter
$ws.Range("C3").Value = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.primary_disease_site IN [''Mouth'']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'

# --- FilesTab row (row 4) ---
# B4: files query, now without the trailing `Study Code` column
$ws.Range("B4").Value = '
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.primary_disease_site IN [''Mouth'']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '''') AS `File Name`, 
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `Format`,
        coalesce(f.file_size, '''') AS `Size`,
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(diag.disease_term,'''') AS Diagnosis'
# C4: StatQuery replaced with the new program/study/case/sample/file counter
$ws.Range("C4").Value = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.primary_disease_site IN [''Mouth'']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'

# Row heights adjusted to fit the revised query text
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# View state: zoom + selection as left by the editing session
$excel.ActiveWindow.Zoom = 70
$ws.Range("C13:C14").Select() | Out-Null
